$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 64, shifting existing rows 64-116 down to 65-117.
$ws.Rows(64).Insert()

# Populate the newly inserted row 64 with a fresh observation. The
# "categorical" columns (Mercado/Region/Producto/... plus Calidad, Unidad de
# comercializacion, Kg/unidad) carry the same values the row above already
# used; the observation-specific columns (Fecha, Volumen, Precios, Origen,
# Precio $/Kg) get the new data point.
$ws.Range("A64").Value = 3
$ws.Range("B64").Value = "Femacal de La Calera"
$ws.Range("C64").Value = "Coquimbo"
$ws.Range("D64").Value = 44484
$ws.Range("E64").Value = 5
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100101
$ws.Range("H64").Value = "Berries"
$ws.Range("I64").Value = 100101001
$ws.Range("J64").Value = "Arándano (blue)"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 45
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 10000
$ws.Range("P64").Value = 10000
$ws.Range("Q64").Value = "`$/bandeja 2 kilos"
$ws.Range("R64").Value = "Provincia de Quillota"
$ws.Range("S64").Value = 5000
$ws.Range("T64").Value = 2
